$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78..186 down to 79..187
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new data record
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 44539
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = 100112009
$ws.Cells.Item(78, 7).Value = "Acelga"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 220
$ws.Cells.Item(78, 11).Value = 600
$ws.Cells.Item(78, 12).Value = 650
$ws.Cells.Item(78, 13).Value = 627
$ws.Cells.Item(78, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(78, 15).Value = "Región de Ñuble"
$ws.Cells.Item(78, 16).Value = 627
$ws.Cells.Item(78, 17).Value = 1
$ws.Cells.Item(78, 18).Value = "Hortaliza"
